$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

$ws.Range("I5").Value = 96.23409753561796
$ws.Range("I6").Value = 96.27694179971012
$ws.Range("I7").Value = 96.26482780896683

$ws.Range("G20").Value = 97.99309799135672
$ws.Range("G21").Value = 98.06967589313628
$ws.Range("G22").Value = 98.05102767438815

$ws.Range("H23").Value = 97.6080679929606
$ws.Range("H24").Value = 97.59876523102058
$ws.Range("H25").Value = 97.57268037750025

$ws.Range("I28").Value = 96.12448977290784
$ws.Range("I29").Value = 96.17446353523682

$ws.Range("G38").Value = 98.02094957675021
$ws.Range("G39").Value = 98.00875161149141

$ws.Range("H40").Value = 97.62350625317688
$ws.Range("H41").Value = 97.59873840950907

$ws.Range("I44").Value = 96.20702259162157
$ws.Range("I45").Value = 96.23701309301885

$ws.Range("G54").Value = 98.0048941774816
$ws.Range("G55").Value = 98.03331257365674

$ws.Range("H56").Value = 97.52481022902361
$ws.Range("H57").Value = 97.45241074553344
